$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The historical price series occupies rows 32-117 (two rows per date: "Primera"
# and "Segunda" quality grades). A new week of data is being recorded, which
# pushes the whole Fecha/Volumen/Origen time series down by one week (2 rows):
#   - rows 34-117 take on the Fecha(D)/Volumen(J)/Origen(O) that used to live
#     two rows above them (old row r -> new row r+2)
#   - rows 32-33 become the newest entry (new date, same J/O as before)
#   - the data that falls off the bottom (old rows 116-117) is appended as new
#     rows 118-119, copied in full
# All other columns (A,B,C,E,F,G,H,I,K,L,M,N,P,Q,R) stay put since they are
# fixed per-row attributes (market, region, product, quality grade, etc.).

$firstRow = 32
$lastRow = 117
$newRowCount = 2

# 1) Snapshot the current D/J/O columns for rows 32-117 before writing anything.
$oldD = @{}
$oldJ = @{}
$oldO = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $oldD[$r] = $ws.Cells.Item($r, 4).Value()
    $oldJ[$r] = $ws.Cells.Item($r, 10).Value()
    $oldO[$r] = $ws.Cells.Item($r, 15).Value()
}

# 2) Snapshot the full rows that will fall off the bottom (old 116-117) so they
#    can be appended as brand-new rows 118-119.
$lastCol = 18
$fallOffRows = @($lastRow - 1, $lastRow)
$fallOffData = @{}
foreach ($r in $fallOffRows) {
    $rowVals = @{}
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $fallOffData[$r] = $rowVals
}

# 3) Shift Fecha/Volumen/Origen down by 2 rows for rows 34-117.
for ($r = $lastRow; $r -ge ($firstRow + $newRowCount); $r--) {
    $srcRow = $r - $newRowCount
    $ws.Cells.Item($r, 4).Value = $oldD[$srcRow]
    $ws.Cells.Item($r, 10).Value = $oldJ[$srcRow]
    $ws.Cells.Item($r, 15).Value = $oldO[$srcRow]
}

# 4) New top entry (rows 32-33): newest reported date, Volumen/Origen unchanged.
$newDate = 44525
$ws.Cells.Item(32, 4).Value = $newDate
$ws.Cells.Item(33, 4).Value = $newDate
$ws.Cells.Item(32, 10).Value = $oldJ[32]
$ws.Cells.Item(33, 10).Value = $oldJ[33]
$ws.Cells.Item(32, 15).Value = $oldO[32]
$ws.Cells.Item(33, 15).Value = $oldO[33]

# 5) Append the rows that fell off the bottom as new rows 118-119 (full copy).
$dateFormat = $ws.Cells.Item($lastRow, 4).NumberFormat()
$destRow = $lastRow + 1
foreach ($r in $fallOffRows) {
    $rowVals = $fallOffData[$r]
    # Set the Fecha (date) number format before the value so the new cell
    # reuses the existing date style instead of allocating a new one.
    $ws.Cells.Item($destRow, 4).NumberFormat = $dateFormat
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $rowVals[$c]
    }
    $destRow++
}
